# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list (rows 16-22 in column E) gets re-ordered from
# descending (2006..1912) to ascending (1912..2006); since every period
# already shares the same "Valor Mora" of 33125 except the one for 2006
# (which was 26500), re-ordering the periods means the 26500 value now
# belongs to the row that ends up holding "1912" instead of "2006".
#
# Net effect on the data grid:
#   F16 (period 2006): 26500 -> 33125
#   F22 (period 1912): 33125 -> 26500

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("F16").Value = 33125
$ws.Range("F22").Value = 26500
